# Update gh-pages output data (广州-漫展信息.xlsx) — refresh "想去人数"
# (want-to-go counts) and flip a handful of sold-out listings from
# "已售罄"/numeric placeholder to "不可售".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 34
$ws1.Range("F3").Value = 1261
$ws1.Range("F4").Value = 12899
$ws1.Range("F5").Value = 736
$ws1.Range("G6").Value = "不可售"
$ws1.Range("G7").Value = "不可售"
$ws1.Range("F8").Value = 59
$ws1.Range("G8").Value = "不可售"
$ws1.Range("F10").Value = 1877
$ws1.Range("G12").Value = "不可售"
$ws1.Range("F18").Value = 298
$ws1.Range("F19").Value = 133
$ws1.Range("F20").Value = 130
$ws1.Range("F22").Value = 220
$ws1.Range("F23").Value = 254
$ws1.Range("F24").Value = 1302
$ws1.Range("F25").Value = 337
$ws1.Range("F27").Value = 103

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 7
$ws2.Range("F5").Value = 4451
$ws2.Range("F6").Value = 163
$ws2.Range("F11").Value = 359
$ws2.Range("F16").Value = 12
$ws2.Range("F17").Value = 13

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 877
$ws3.Range("F3").Value = 4228

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types — aggregate of the sheets above)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 877
$ws4.Range("F3").Value = 34
$ws4.Range("F5").Value = 7
$ws4.Range("F6").Value = 1261
$ws4.Range("F7").Value = 12899
$ws4.Range("F9").Value = 736
$ws4.Range("F10").Value = 4228
$ws4.Range("G11").Value = "不可售"
$ws4.Range("G12").Value = "不可售"
$ws4.Range("F13").Value = 59
$ws4.Range("G13").Value = "不可售"
$ws4.Range("F15").Value = 1877
$ws4.Range("G17").Value = "不可售"
$ws4.Range("F19").Value = 4451
$ws4.Range("F21").Value = 163
$ws4.Range("F22").Value = 163
$ws4.Range("F28").Value = 359
$ws4.Range("F32").Value = 298
$ws4.Range("F33").Value = 133
$ws4.Range("F34").Value = 130
$ws4.Range("F37").Value = 220
$ws4.Range("F40").Value = 254
$ws4.Range("F41").Value = 1302
$ws4.Range("F42").Value = 12
$ws4.Range("F43").Value = 337
$ws4.Range("F45").Value = 103
$ws4.Range("F46").Value = 13
